# Backup QR Scanner data - 2025-12-08T18:47:27.028Z - Cache Bust: 1765219647028
#
# 1. Rename the worksheet "Scanner" -> "Session".
# 2. Normalize the "User" column (F) e-mail addresses to lower-case
#    ("Emp16...." -> "emp16....").
# 3. Remove the last three scanner log rows (13-15); the scanner backup now
#    only covers the first 11 log entries (rows 2-12 of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename the sheet ---------------------------------------------------
$ws.Name = "Session"

# --- 2. Lower-case the "User" e-mail for every remaining data row ----------
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 6).Value = "emp16.farida.m.abdelaziz@gmail.com"
}

# --- 3. Drop the trailing rows (student IDs 191076 / 190795 / 180752) ------
$ws.Rows("13:15").Delete()

# --- 4. Keep the "number stored as text" ignored-error hint in sync with
#        the shrunk data range (A1:F15 -> A1:F12). --------------------------
$errRange = $ws.Range("A1:F12")
$numberAsText = $errRange.Errors.Item(3)
$numberAsText.Ignore = $true
